# Updates the cryptos price list (Price / Volume(1h) columns, plus the
# PEPE <-> Binance-PegBSC-USD row swap) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes $Value into $Cell as literal text. Some "Price" values look
    # like plain numbers (e.g. "596.04"), and Excel would otherwise auto-
    # convert them to numeric cells; prefixing with an apostrophe forces
    # text, then resetting the style drops the quote-prefix formatting
    # that COM tags on so the cell keeps its original (unstyled) look.
    param($Cell, [string]$Value)
    $Cell.Value = "'" + $Value
    $Cell.Style = "Normal"
}


$ws.Range("D2").Value = '67.307.80'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '2.619.03'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("E4").Value = '  -0.10%  '

Set-TextValue $ws.Range("D5") '596.04'
$ws.Range("E5").Value = '  +0.19%  '

Set-TextValue $ws.Range("D6") '153.01'
$ws.Range("E6").Value = '  -1.59%  '

$ws.Range("E7").Value = '  +0.03%  '

Set-TextValue $ws.Range("D8") '0.557'
$ws.Range("E8").Value = '  +2.64%  '

$ws.Range("D9").Value = '2.617.48'
$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  -3.88%  '

$ws.Range("E12").Value = '  -0.93%  '

$ws.Range("E13").Value = '  -2.34%  '

Set-TextValue $ws.Range("D14") '27.74'
$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").Value = '3.095.40'
$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("E16").Value = '  -4.79%  '

$ws.Range("D17").Value = '67.159.96'
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").Value = '2.616.95'
$ws.Range("E18").Value = '  -0.02%  '

Set-TextValue $ws.Range("D19") '11.10'
$ws.Range("E19").Value = '  -2.17%  '

Set-TextValue $ws.Range("D20") '363.54'
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("E21").Value = '  -4.34%  '

$ws.Range("E22").Value = '  -0.58%  '

$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("E24").Value = '  -0.15%  '

Set-TextValue $ws.Range("D25") '71.09'
$ws.Range("E25").Value = '  +4.93%  '

Set-TextValue $ws.Range("D26") '10.03'
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("D27").Value = '2.743.50'

Set-TextValue $ws.Range("D28") '586.06'
$ws.Range("E28").Value = '  -3.78%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D29") '0.0000102'
$ws.Range("E29").Value = '  -2.68%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("E31").Value = '  -4.23%  '

$ws.Range("E32").Value = '  -2.28%  '

$ws.Range("E33").Value = '  -1.56%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  -5.71%  '

$ws.Range("E36").Value = '  -3.00%  '

$ws.Range("E37").Value = '  -2.35%  '

Set-TextValue $ws.Range("D38") '157.39'
$ws.Range("E38").Value = '  +1.55%  '

Set-TextValue $ws.Range("D39") '19.09'
$ws.Range("E39").Value = '  -2.84%  '

$ws.Range("E40").Value = '  -0.63%  '

$ws.Range("E41").Value = '  -3.70%  '

$ws.Range("E42").Value = '  -1.79%  '

Set-TextValue $ws.Range("D43") '2.57'
$ws.Range("E43").Value = '  -2.42%  '

Set-TextValue $ws.Range("D44") '41.14'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("E45").Value = '  -0.03%  '

Set-TextValue $ws.Range("D46") '16.37'
$ws.Range("E46").Value = '  -0.63%  '

Set-TextValue $ws.Range("D47") '156.74'
$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("E49").Value = '  -0.87%  '

Set-TextValue $ws.Range("D50") '21.95'

Set-TextValue $ws.Range("D51") '0.623'
$ws.Range("E51").Value = '  -0.95%  '
